$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.5
$ws.Range("J2").Value = 3.3
$ws.Range("N2").Value = 3.55
$ws.Range("AA2").Value = 48
$ws.Range("H3").Value = 4.4
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 4.2
$ws.Range("K3").Value = 4.3
$ws.Range("P3").Value = 2.46
$ws.Range("R3").Value = 1.59
$ws.Range("T3").Value = 1.64
$ws.Range("U3").Value = 2.46
$ws.Range("V3").Value = 1.28
$ws.Range("W3").Value = 2.16
$ws.Range("Y3").Value = 22
$ws.Range("AD3").Value = 18
$ws.Range("AJ3").Value = 21
$ws.Range("AM3").Value = 70
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 12.5
$ws.Range("H4").Value = 1.31
$ws.Range("J4").Value = 5.6
$ws.Range("N4").Value = 4.9
$ws.Range("P4").Value = 2.36
$ws.Range("Q4").Value = 1.6
$ws.Range("T4").Value = 1.98
$ws.Range("AD4").Value = 11
$ws.Range("AG4").Value = 46
$ws.Range("AJ4").Value = 450
$ws.Range("AL4").Value = 150
$ws.Range("AN4").Value = 250
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4.6
$ws.Range("P5").Value = 2.3
$ws.Range("Q5").Value = 1.62
$ws.Range("S5").Value = 2.28
$ws.Range("T5").Value = 1.66
$ws.Range("AF5").Value = 20
$ws.Range("F6").Value = 1.37
$ws.Range("J6").Value = 5.9
$ws.Range("T6").Value = 1.78
$ws.Range("U6").Value = 2.2
$ws.Range("H7").Value = 2.12
$ws.Range("I7").Value = 2.16
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 3.9
$ws.Range("P7").Value = 2.4
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 1.54
$ws.Range("S7").Value = 2.74
$ws.Range("U7").Value = 2.5
$ws.Range("AB7").Value = 17.5
$ws.Range("AJ7").Value = 65
$ws.Range("AN7").Value = 27
$ws.Range("F8").Value = 1.69
$ws.Range("G8").Value = 1.7
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 6.4
$ws.Range("L8").Value = 1.39
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 1.97
$ws.Range("Q8").Value = 2
$ws.Range("S8").Value = 3.55
$ws.Range("U8").Value = 1.97
$ws.Range("V8").Value = 1.18
$ws.Range("W8").Value = 2.42
$ws.Range("AB8").Value = 8
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AG8").Value = 9.6
$ws.Range("AM8").Value = 130
$ws.Range("F9").Value = 1.09
$ws.Range("I9").Value = 970
$ws.Range("J9").Value = 1.09
$ws.Range("R12").Value = 1.33
$ws.Range("F13").Value = 2.6
$ws.Range("I13").Value = 2.78
$ws.Range("J13").Value = 3.85
$ws.Range("K13").Value = 4.2
$ws.Range("Q13").Value = 1.61
$ws.Range("R13").Value = 1.56
$ws.Range("V13").Value = 1.56
$ws.Range("W13").Value = 1.58
$ws.Range("F14").Value = 1.63
$ws.Range("G14").Value = 1.66
$ws.Range("H14").Value = 5.2
$ws.Range("T14").Value = 1.71
$ws.Range("V14").Value = 1.21
$ws.Range("W14").Value = 2.5
$ws.Range("X14").Value = 23
$ws.Range("AE14").Value = 65
$ws.Range("AK14").Value = 17
$ws.Range("AN14").Value = 8
